$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 6.022537231445312
$ws.Range("B1").Value = 4.748425483703613
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 4.369626045227051
$ws.Range("E1").Value = 2.588987350463867
